$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(4).NumberFormat = "@"

$ws.Range('D2').Value = '27.068.06'
$ws.Range('E2').Value = '  -2.49%  '
$ws.Range('D3').Value = '1.865.61'
$ws.Range('E3').Value = '  -2.17%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('D5').Value = '306.39'
$ws.Range('E5').Value = '  -1.99%  '
$ws.Range('D7').Value = '0.5128'
$ws.Range('E7').Value = '  -1.91%  '
$ws.Range('D8').Value = '0.3754'
$ws.Range('E8').Value = '  -0.84%  '
$ws.Range('E9').Value = '  -1.17%  '
$ws.Range('D10').Value = '0.8894'
$ws.Range('E10').Value = '  -1.43%  '
$ws.Range('D11').Value = '20.69'
$ws.Range('E11').Value = '  -3.16%  '
$ws.Range('D12').Value = '0.07586'
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('D13').Value = '1.841.64'
$ws.Range('E13').Value = '  -3.07%  '
$ws.Range('D14').Value = '5.309'
$ws.Range('E14').Value = '  -2.67%  '
$ws.Range('D15').Value = '89.47'
$ws.Range('E15').Value = '  -2.92%  '
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').Value = '0.000008452'
$ws.Range('E17').Value = '  -2.90%  '
$ws.Range('E18').Value = '  -2.73%  '
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').Value = '27.113.32'
$ws.Range('E20').Value = '  -2.40%  '
$ws.Range('D21').Value = '5.042'
$ws.Range('E21').Value = '  -2.05%  '
$ws.Range('D22').Value = '2.095.96'
$ws.Range('E22').Value = '  -0.93%  '
$ws.Range('D23').Value = '10.53'
$ws.Range('E23').Value = '  -2.82%  '
$ws.Range('D24').Value = '6.455'
$ws.Range('E24').Value = '  -1.95%  '
$ws.Range('D25').Value = '1.843'
$ws.Range('E25').Value = '  -1.62%  '
$ws.Range('D26').Value = '147.52'
$ws.Range('E26').Value = '  -3.76%  '
$ws.Range('D27').Value = '17.97'
$ws.Range('E27').Value = '  -1.94%  '
$ws.Range('D28').Value = '2.117'
$ws.Range('E28').Value = '  -2.49%  '
$ws.Range('D29').Value = '112.82'
$ws.Range('E29').Value = '  -1.57%  '
$ws.Range('D30').Value = '4.662'
$ws.Range('E30').Value = '  -4.22%  '
$ws.Range('D31').Value = '4.709'
$ws.Range('E31').Value = '  -3.33%  '
$ws.Range('D32').Value = '0.09106'
$ws.Range('E32').Value = '  +1.35%  '
$ws.Range('D34').Value = '3.061'
$ws.Range('E34').Value = '  -3.56%  '
$ws.Range('D35').Value = '1.156'
$ws.Range('E35').Value = '  -6.14%  '
$ws.Range('D36').Value = '0.7279'
$ws.Range('E36').Value = '  -6.03%  '
$ws.Range('D37').Value = '0.02042'
$ws.Range('E37').Value = '  -2.10%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '3.046'
$ws.Range('E38').Value = '  -0.91%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '2.476'
$ws.Range('E39').Value = '  -5.95%  '
$ws.Range('D40').Value = '1.075'
$ws.Range('E40').Value = '  -1.80%  '
$ws.Range('D41').Value = '0.5334'
$ws.Range('E41').Value = '  -3.40%  '
$ws.Range('D42').Value = '6.570'
$ws.Range('E42').Value = '  -1.55%  '
$ws.Range('D43').Value = '117.44'
$ws.Range('E43').Value = '  +2.55%  '
$ws.Range('D44').Value = '8.277'
$ws.Range('E44').Value = '  -2.95%  '
$ws.Range('D46').Value = '0.4634'
$ws.Range('E46').Value = '  -3.78%  '
$ws.Range('E47').Value = '  +0.20%  '
$ws.Range('E48').Value = '  -4.49%  '
$ws.Range('D49').Value = '1.575'
$ws.Range('E49').Value = '  -2.70%  '
$ws.Range('D50').Value = '36.60'
$ws.Range('E50').Value = '  -0.74%  '
$ws.Range('E51').Value = '  -4.22%  '
